$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("K" header) values recomputed from Strike# to K; updating per-row.
$gValues = @(3,1,1,2,0,2,0,0,1,1,1,1,1,1,1,1,1,1,1,0,1,0,4,1,0,1,0,0,1,2,1,2,2,1,0,1,0,0,0,0,0,1,0,2,0,1,1,0,0,0,3,1,1,2,3,2,1,4,2,3,3,1,0,1,2,1,2,1)

$startRow = 2
for ($i = 0; $i -lt $gValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $gValues[$i]
}

Write-Host "Updated G2:G69 with new K values"
